$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.768105626106262
$ws.Range("B1").Value = 3.437867403030396
$ws.Range("C1").Value = 3.849934816360474
$ws.Range("D1").Value = 3.592831373214722
$ws.Range("E1").Value = 1.033917307853699
